$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.2708946666666667
$ws.Cells.Item(2, 8).Value = 0.812684
$ws.Cells.Item(2, 9).Value = 0.1616296696421007
$ws.Cells.Item(2, 10).Value = 0.1616296696421007
$ws.Cells.Item(2, 13).Value = 15.50220733333333
$ws.Cells.Item(2, 14).Value = 46.506622
$ws.Cells.Item(2, 15).Value = 0.5994675913188158
$ws.Cells.Item(2, 16).Value = 0.5994675913188158
$ws.Cells.Item(2, 17).Value = 4.199465288160889
$ws.Cells.Item(2, 18).Value = 37.795187593448
$ws.Cells.Item(2, 19).Value = 0.09689174874600599
$ws.Cells.Item(2, 20).Value = 0.09689174874600599

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.2708946666666667
$ws.Cells.Item(3, 8).Value = 0.812684
$ws.Cells.Item(3, 9).Value = 0.1616296696421007
$ws.Cells.Item(3, 10).Value = 0.1616296696421007
$ws.Cells.Item(3, 15).Value = 0.04399860030713892
$ws.Cells.Item(3, 16).Value = 0.04399860030713892
$ws.Cells.Item(3, 17).Value = 0.3082244935226667
$ws.Cells.Item(3, 18).Value = 2.774020441704
$ws.Cells.Item(3, 19).Value = 0.007111479232357693
$ws.Cells.Item(3, 20).Value = 0.007111479232357693

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2708946666666667
$ws.Cells.Item(4, 8).Value = 0.812684
$ws.Cells.Item(4, 9).Value = 0.1616296696421007
$ws.Cells.Item(4, 10).Value = 0.1616296696421007
$ws.Cells.Item(4, 13).Value = 8.848210666666667
$ws.Cells.Item(4, 14).Value = 26.544632
$ws.Cells.Item(4, 15).Value = 0.3421587275782868
$ws.Cells.Item(4, 16).Value = 0.3421587275782868
$ws.Cells.Item(4, 17).Value = 2.396933079143111
$ws.Cells.Item(4, 18).Value = 21.572397712288
$ws.Cells.Item(4, 19).Value = 0.05530300210364001
$ws.Cells.Item(4, 20).Value = 0.05530300210364001

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2708946666666667
$ws.Cells.Item(5, 8).Value = 0.812684
$ws.Cells.Item(5, 9).Value = 0.1616296696421007
$ws.Cells.Item(5, 10).Value = 0.1616296696421007
$ws.Cells.Item(5, 13).Value = 0.371739
$ws.Cells.Item(5, 14).Value = 1.115217
$ws.Cells.Item(5, 15).Value = 0.01437508079575842
$ws.Cells.Item(5, 16).Value = 0.01437508079575841
$ws.Cells.Item(5, 17).Value = 0.100702112492
$ws.Cells.Item(5, 18).Value = 0.9063190124279998
$ws.Cells.Item(5, 19).Value = 0.002323439560096938
$ws.Cells.Item(5, 20).Value = 0.002323439560096938

# Row 6
$ws.Cells.Item(6, 9).Value = 0.6313295261673385
$ws.Cells.Item(6, 10).Value = 0.6313295261673384
$ws.Cells.Item(6, 13).Value = 15.50220733333333
$ws.Cells.Item(6, 14).Value = 46.506622
$ws.Cells.Item(6, 15).Value = 0.5994675913188158
$ws.Cells.Item(6, 16).Value = 0.5994675913188158
$ws.Cells.Item(6, 17).Value = 16.40321629315645
$ws.Cells.Item(6, 18).Value = 147.628946638408
$ws.Cells.Item(6, 19).Value = 0.3784615903799837
$ws.Cells.Item(6, 20).Value = 0.3784615903799836

# Row 7
$ws.Cells.Item(7, 9).Value = 0.6313295261673385
$ws.Cells.Item(7, 10).Value = 0.6313295261673384
$ws.Cells.Item(7, 15).Value = 0.04399860030713892
$ws.Cells.Item(7, 16).Value = 0.04399860030713892
$ws.Cells.Item(7, 19).Value = 0.02777761548393213
$ws.Cells.Item(7, 20).Value = 0.02777761548393213

# Row 8
$ws.Cells.Item(8, 9).Value = 0.6313295261673385
$ws.Cells.Item(8, 10).Value = 0.6313295261673384
$ws.Cells.Item(8, 13).Value = 8.848210666666667
$ws.Cells.Item(8, 14).Value = 26.544632
$ws.Cells.Item(8, 15).Value = 0.3421587275782868
$ws.Cells.Item(8, 16).Value = 0.3421587275782868
$ws.Cells.Item(8, 17).Value = 9.362480468227556
$ws.Cells.Item(8, 18).Value = 84.26232421404801
$ws.Cells.Item(8, 19).Value = 0.2160149073560193
$ws.Cells.Item(8, 20).Value = 0.2160149073560192

# Row 9
$ws.Cells.Item(9, 9).Value = 0.6313295261673385
$ws.Cells.Item(9, 10).Value = 0.6313295261673384
$ws.Cells.Item(9, 13).Value = 0.371739
$ws.Cells.Item(9, 14).Value = 1.115217
$ws.Cells.Item(9, 15).Value = 0.01437508079575842
$ws.Cells.Item(9, 16).Value = 0.01437508079575841
$ws.Cells.Item(9, 17).Value = 0.3933449663320001
$ws.Cells.Item(9, 18).Value = 3.540104696988
$ws.Cells.Item(9, 19).Value = 0.009075412947403367
$ws.Cells.Item(9, 20).Value = 0.009075412947403364

# Row 10
$ws.Cells.Item(10, 9).Value = 0.2070408041905609
$ws.Cells.Item(10, 10).Value = 0.2070408041905609
$ws.Cells.Item(10, 13).Value = 15.50220733333333
$ws.Cells.Item(10, 14).Value = 46.506622
$ws.Cells.Item(10, 15).Value = 0.5994675913188158
$ws.Cells.Item(10, 16).Value = 0.5994675913188158
$ws.Cells.Item(10, 17).Value = 5.379338288300889
$ws.Cells.Item(10, 18).Value = 48.414044594708
$ws.Cells.Item(10, 19).Value = 0.1241142521928261
$ws.Cells.Item(10, 20).Value = 0.1241142521928261

# Row 11
$ws.Cells.Item(11, 9).Value = 0.2070408041905609
$ws.Cells.Item(11, 10).Value = 0.2070408041905609
$ws.Cells.Item(11, 15).Value = 0.04399860030713892
$ws.Cells.Item(11, 16).Value = 0.04399860030713892
$ws.Cells.Item(11, 19).Value = 0.009109505590849103
$ws.Cells.Item(11, 20).Value = 0.009109505590849103

# Row 12
$ws.Cells.Item(12, 9).Value = 0.2070408041905609
$ws.Cells.Item(12, 10).Value = 0.2070408041905609
$ws.Cells.Item(12, 13).Value = 8.848210666666667
$ws.Cells.Item(12, 14).Value = 26.544632
$ws.Cells.Item(12, 15).Value = 0.3421587275782868
$ws.Cells.Item(12, 16).Value = 0.3421587275782868
$ws.Cells.Item(12, 17).Value = 3.070370392983111
$ws.Cells.Item(12, 18).Value = 27.633333536848
$ws.Cells.Item(12, 19).Value = 0.07084081811862755
$ws.Cells.Item(12, 20).Value = 0.07084081811862755

# Row 13
$ws.Cells.Item(13, 9).Value = 0.2070408041905609
$ws.Cells.Item(13, 10).Value = 0.2070408041905609
$ws.Cells.Item(13, 13).Value = 0.371739
$ws.Cells.Item(13, 14).Value = 1.115217
$ws.Cells.Item(13, 15).Value = 0.01437508079575842
$ws.Cells.Item(13, 16).Value = 0.01437508079575841
$ws.Cells.Item(13, 17).Value = 0.128995167782
$ws.Cells.Item(13, 18).Value = 1.160956510038
$ws.Cells.Item(13, 19).Value = 0.002976228288258111
$ws.Cells.Item(13, 20).Value = 0.00297622828825811

